# Practica2 "tabla Tokens" - add the missing DOSPUN (":") token row at the
# end of the tokens table, mirroring the formatting of the preceding rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "DOSPUN"
$newRow.Cells.Item(2).Range.Text = [char]8220 + ":" + [char]8221
$newRow.Cells.Item(3).Range.Text = "291"
$newRow.Cells.Item(4).Range.Text = "{:}"
